$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Agency (column D) names.
# Most rows belong to the "Secretaría de Seguridad Pública del Estado de Tabasco",
# but row 11 (Rafael Collado Gómez) is actually a different, similarly-named
# agency: "Secretaría de Seguridad Pública de Tabasco" (apparent duplication
# under investigation per the Agency Completion Search).

$ws.Range("D1:D18").Value = "Secretaría de Seguridad Pública del Estado de Tabasco"
$ws.Range("D11").Value = "Secretaría de Seguridad Pública de Tabasco"
